# New stepsizes -> edited stepsize 2 and 5 to use the absolute value of
# gamma (N<row>) in the denominator instead of the quadratic term
# (Column3 / Column5 / Column7) for every data row of the Lower_bounds
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column "Column3" relative-improvement (cells E5:E20) ------------
# Row 5 used to hold a hard-coded constant (-1E-3); rows 6:20 used a
# formula whose numerator had a hard-coded literal and whose
# denominator referenced the (self) Column3 value. All of them now
# read "(Column3 - N<row>) / N<row>".
for ($r = 5; $r -le 20; $r++) {
    $ws.Range("E$r").Formula = "=(Lower_bounds[[#This Row],[Column3]]-N$r)/N$r"
}

# --- Column "Column5" relative-improvement (cells H5:H8) --------------
# Only rows 5-8 carry a formula in column H (the rest hold the " "
# placeholder string); the denominator switches from the Column5
# self-reference to N<row>.
for ($r = 5; $r -le 8; $r++) {
    $ws.Range("H$r").Formula = "=(Lower_bounds[[#This Row],[Column5]]-N$r)/N$r"
}

# --- Column "Column7" relative-improvement (cells K5:K20) --------------
# Denominator switches from the Column7 self-reference to N<row>.
for ($r = 5; $r -le 20; $r++) {
    $ws.Range("K$r").Formula = "=(Lower_bounds[[#This Row],[Column7]]-N$r)/N$r"
}

# --- New summary labels / averages in column M ------------------------
$ws.Range("M7").Value = "IR av:"
$ws.Range("M9").Value = "AP av:"
$ws.Range("M12").Value = "AV IR 8:"

$ws.Range("M8").Formula = "=AVERAGE(E5:E20)"
$ws.Range("M8").Style = "Prozent"

$ws.Range("M10").Formula = "=AVERAGE(H5:H8)"
$ws.Range("M10").Style = "Prozent"

$ws.Range("M13").Formula = "=AVERAGE(E13:E20)"
$ws.Range("M13").Style = "Prozent"

# --- Sheet view: selection moves from M7 to M14 ------------------------
$ws.Range("M14").Select()
